$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Category (column F) values per row, as specified by the diff.
# Mapping pattern: Financing -> Investing, Investing -> Operating, Operating -> Investing
$updates = @{
    3 = "Investing Activities (XGBoost)"
    4 = "Operating Activities (XGBoost)"
    5 = "Investing Activities (XGBoost)"
    6 = "Operating Activities (XGBoost)"
    10 = "Investing Activities (XGBoost)"
    11 = "Operating Activities (XGBoost)"
    16 = "Investing Activities (XGBoost)"
    17 = "Investing Activities (XGBoost)"
    22 = "Operating Activities (XGBoost)"
    27 = "Investing Activities (XGBoost)"
    29 = "Operating Activities (XGBoost)"
    32 = "Investing Activities (XGBoost)"
    33 = "Operating Activities (XGBoost)"
    34 = "Operating Activities (XGBoost)"
    35 = "Operating Activities (XGBoost)"
    37 = "Operating Activities (XGBoost)"
    47 = "Operating Activities (XGBoost)"
    49 = "Operating Activities (XGBoost)"
    51 = "Operating Activities (XGBoost)"
    53 = "Investing Activities (XGBoost)"
    58 = "Operating Activities (XGBoost)"
    60 = "Operating Activities (XGBoost)"
    64 = "Investing Activities (XGBoost)"
    65 = "Operating Activities (XGBoost)"
    67 = "Investing Activities (XGBoost)"
    68 = "Investing Activities (XGBoost)"
    77 = "Operating Activities (XGBoost)"
    81 = "Operating Activities (XGBoost)"
    86 = "Operating Activities (XGBoost)"
    87 = "Operating Activities (XGBoost)"
    91 = "Operating Activities (XGBoost)"
    92 = "Operating Activities (XGBoost)"
    100 = "Investing Activities (XGBoost)"
    104 = "Investing Activities (XGBoost)"
    105 = "Operating Activities (XGBoost)"
    107 = "Operating Activities (XGBoost)"
    111 = "Operating Activities (XGBoost)"
    115 = "Operating Activities (XGBoost)"
    117 = "Operating Activities (XGBoost)"
    119 = "Investing Activities (XGBoost)"
    122 = "Investing Activities (XGBoost)"
    126 = "Operating Activities (XGBoost)"
    128 = "Investing Activities (XGBoost)"
    130 = "Operating Activities (XGBoost)"
    133 = "Operating Activities (XGBoost)"
    134 = "Operating Activities (XGBoost)"
    135 = "Investing Activities (XGBoost)"
    143 = "Investing Activities (XGBoost)"
    144 = "Operating Activities (XGBoost)"
    145 = "Operating Activities (XGBoost)"
    147 = "Operating Activities (XGBoost)"
    155 = "Operating Activities (XGBoost)"
    157 = "Investing Activities (XGBoost)"
    160 = "Investing Activities (XGBoost)"
    162 = "Investing Activities (XGBoost)"
    170 = "Operating Activities (XGBoost)"
    171 = "Operating Activities (XGBoost)"
    176 = "Operating Activities (XGBoost)"
    178 = "Operating Activities (XGBoost)"
    180 = "Operating Activities (XGBoost)"
    185 = "Operating Activities (XGBoost)"
    189 = "Operating Activities (XGBoost)"
    193 = "Investing Activities (XGBoost)"
    194 = "Investing Activities (XGBoost)"
    200 = "Investing Activities (XGBoost)"
    206 = "Operating Activities (XGBoost)"
    207 = "Operating Activities (XGBoost)"
    216 = "Operating Activities (XGBoost)"
    217 = "Investing Activities (XGBoost)"
    223 = "Operating Activities (XGBoost)"
    224 = "Operating Activities (XGBoost)"
    228 = "Operating Activities (XGBoost)"
    239 = "Investing Activities (XGBoost)"
    241 = "Operating Activities (XGBoost)"
    246 = "Operating Activities (XGBoost)"
    247 = "Operating Activities (XGBoost)"
    249 = "Investing Activities (XGBoost)"
    257 = "Investing Activities (XGBoost)"
    258 = "Operating Activities (XGBoost)"
    259 = "Investing Activities (XGBoost)"
    262 = "Operating Activities (XGBoost)"
    265 = "Investing Activities (XGBoost)"
    266 = "Investing Activities (XGBoost)"
    267 = "Investing Activities (XGBoost)"
    268 = "Operating Activities (XGBoost)"
    273 = "Operating Activities (XGBoost)"
    276 = "Investing Activities (XGBoost)"
    282 = "Operating Activities (XGBoost)"
    287 = "Investing Activities (XGBoost)"
    292 = "Operating Activities (XGBoost)"
    294 = "Operating Activities (XGBoost)"
    295 = "Investing Activities (XGBoost)"
    298 = "Operating Activities (XGBoost)"
    299 = "Operating Activities (XGBoost)"
    300 = "Operating Activities (XGBoost)"
    301 = "Investing Activities (XGBoost)"
    304 = "Operating Activities (XGBoost)"
    319 = "Operating Activities (XGBoost)"
    320 = "Operating Activities (XGBoost)"
    327 = "Operating Activities (XGBoost)"
    328 = "Investing Activities (XGBoost)"
    330 = "Investing Activities (XGBoost)"
    339 = "Operating Activities (XGBoost)"
    354 = "Investing Activities (XGBoost)"
    369 = "Investing Activities (XGBoost)"
    375 = "Investing Activities (XGBoost)"
    377 = "Investing Activities (XGBoost)"
    379 = "Investing Activities (XGBoost)"
    385 = "Operating Activities (XGBoost)"
    386 = "Operating Activities (XGBoost)"
    392 = "Operating Activities (XGBoost)"
    399 = "Operating Activities (XGBoost)"
    400 = "Operating Activities (XGBoost)"
    403 = "Operating Activities (XGBoost)"
    404 = "Operating Activities (XGBoost)"
    405 = "Operating Activities (XGBoost)"
    408 = "Operating Activities (XGBoost)"
    411 = "Operating Activities (XGBoost)"
    412 = "Operating Activities (XGBoost)"
    415 = "Investing Activities (XGBoost)"
    419 = "Operating Activities (XGBoost)"
    430 = "Operating Activities (XGBoost)"
    433 = "Operating Activities (XGBoost)"
    434 = "Operating Activities (XGBoost)"
    435 = "Operating Activities (XGBoost)"
    437 = "Operating Activities (XGBoost)"
    438 = "Operating Activities (XGBoost)"
    442 = "Investing Activities (XGBoost)"
    446 = "Investing Activities (XGBoost)"
    447 = "Operating Activities (XGBoost)"
    449 = "Investing Activities (XGBoost)"
    450 = "Investing Activities (XGBoost)"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
